$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: classical-best-embeddings vs. classical-best-tfidf -> classical-best-embed vs. classical-best-tfidf
$ws.Range("A2").Value = "classical-best-embed vs. classical-best-tfidf"
$ws.Range("C2").Value = 0.096
$ws.Range("D2").Value = 0.06
$ws.Range("I2").Value = 0.079
$ws.Range("J2").Value = 0.084

# Row 3: BERT-base vs. classical-best-tfidf
$ws.Range("C3").Value = 0.091
$ws.Range("D3").Value = 0.164
$ws.Range("E3").Value = 0.17
$ws.Range("F3").Value = 0.17
$ws.Range("G3").Value = 0.186
$ws.Range("H3").Value = 0.217
$ws.Range("I3").Value = 0.149
$ws.Range("J3").Value = 0.166

# Row 4: BERT-base vs. classical-best-embeddings -> BERT-base vs. classical-best-embed
$ws.Range("A4").Value = "BERT-base vs. classical-best-embed"
$ws.Range("D4").Value = 0.104
$ws.Range("E4").Value = 0.086
$ws.Range("F4").Value = 0.094
$ws.Range("G4").Value = 0.115
$ws.Range("H4").Value = 0.102
$ws.Range("I4").Value = 0.07
$ws.Range("J4").Value = 0.083

# Row 5: BERT-base-nli vs. classical-best-tfidf
$ws.Range("B5").Value = 0.359
$ws.Range("C5").Value = 0.284
$ws.Range("D5").Value = 0.248
$ws.Range("E5").Value = 0.247
$ws.Range("F5").Value = 0.236
$ws.Range("G5").Value = 0.225
$ws.Range("H5").Value = 0.247
$ws.Range("I5").Value = 0.254
$ws.Range("J5").Value = 0.248

# Row 6: BERT-base-nli vs. classical-best-embeddings -> BERT-base-nli vs. classical-best-embed
$ws.Range("A6").Value = "BERT-base-nli vs. classical-best-embed"
$ws.Range("B6").Value = 0.359
$ws.Range("C6").Value = 0.188
$ws.Range("D6").Value = 0.188
$ws.Range("E6").Value = 0.163
$ws.Range("F6").Value = 0.16
$ws.Range("G6").Value = 0.154
$ws.Range("H6").Value = 0.132
$ws.Range("I6").Value = 0.175
$ws.Range("J6").Value = 0.164

# Row 7: BERT-base-nli vs. BERT-base
$ws.Range("B7").Value = 0.359
$ws.Range("C7").Value = 0.193
$ws.Range("D7").Value = 0.084
$ws.Range("E7").Value = 0.077
$ws.Range("I7").Value = 0.105
$ws.Range("J7").Value = 0.082
